$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.797.33'
$ws.Range("E2").Value = '  -1.50%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.749.91'
$ws.Range("E3").Value = '  -3.91%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.41%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.93'
$ws.Range("E5").Value = '  -2.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4222'
$ws.Range("E7").Value = '  -4.58%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3621'
$ws.Range("E8").Value = '  -3.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.42'
$ws.Range("E9").Value = '  -5.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07424'
$ws.Range("E10").Value = '  -4.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.086'
$ws.Range("E11").Value = '  -3.11%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.71'
$ws.Range("E13").Value = '  -6.19%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.034'
$ws.Range("E14").Value = '  -4.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.264'
$ws.Range("E15").Value = '  -3.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.776.32'
$ws.Range("E16").Value = '  -3.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '90.65'
$ws.Range("E17").Value = '  -2.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001051'
$ws.Range("E18").Value = '  -2.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06341'
$ws.Range("E19").Value = '  -2.21%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.00'
$ws.Range("E21").Value = '  -3.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.907'
$ws.Range("E22").Value = '  -6.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.809.77'
$ws.Range("E23").Value = '  -1.75%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.12'
$ws.Range("E24").Value = '  -4.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.097'
$ws.Range("E25").Value = '  -2.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '156.80'
$ws.Range("E26").Value = '  +0.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.09'
$ws.Range("E27").Value = '  -2.90%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.961.19'
$ws.Range("E28").Value = '  -3.64%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.121'
$ws.Range("E29").Value = '  -9.48%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.55'
$ws.Range("E30").Value = '  -3.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.120'
$ws.Range("E31").Value = '  -6.32%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.637'
$ws.Range("E32").Value = '  -0.72%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.543'
$ws.Range("E33").Value = '  -5.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08828'
$ws.Range("E34").Value = '  -4.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.25'
$ws.Range("E35").Value = '  -5.78%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02275'
$ws.Range("E36").Value = '  -2.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2091'
$ws.Range("E37").Value = '  -3.97%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06008'
$ws.Range("E38").Value = '  -2.72%  '
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6293'
$ws.Range("E39").Value = '  -4.22%  '
$ws.Range("B40").Value = 'InternetComputer(DFINITY)'
$ws.Range("C40").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.925'
$ws.Range("E40").Value = '  -4.64%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.170'
$ws.Range("E41").Value = '  -3.08%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.000'
$ws.Range("E42").Value = '  -0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.397'
$ws.Range("E43").Value = '  +0.40%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.840'
$ws.Range("E44").Value = '  -3.32%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.31'
$ws.Range("E45").Value = '  -4.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5853'
$ws.Range("E46").Value = '  -3.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.674'
$ws.Range("E47").Value = '  -2.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.48'
$ws.Range("E48").Value = '  -3.34%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.965'
$ws.Range("E49").Value = '  -3.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.173'
$ws.Range("E50").Value = '  +1.79%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06793'
$ws.Range("E51").Value = '  -2.98%  '
